$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.446.57'
$ws.Range("E2").Value = '  +6.54%  '
$ws.Range("D3").Value = '1.814.29'
$ws.Range("E3").Value = '  +6.30%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '344.74'
$ws.Range("D5").Style = $ws.Range("C5").Style
$ws.Range("E5").Value = '  +4.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9999'
$ws.Range("D6").Style = $ws.Range("C6").Style
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3845'
$ws.Range("D7").Style = $ws.Range("C7").Style
$ws.Range("E7").Value = '  +4.44%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '50.25'
$ws.Range("D8").Style = $ws.Range("C8").Style
$ws.Range("E8").Value = '  +3.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3528'
$ws.Range("D9").Style = $ws.Range("C9").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.238'
$ws.Range("D10").Style = $ws.Range("C10").Style
$ws.Range("E10").Value = '  +5.83%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07770'
$ws.Range("D11").Style = $ws.Range("C11").Style
$ws.Range("E11").Value = '  +5.68%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.000'
$ws.Range("D12").Style = $ws.Range("C12").Style
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("E13").Value = '  +13.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.633'
$ws.Range("D14").Style = $ws.Range("C14").Style
$ws.Range("E14").Value = '  +6.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.237'
$ws.Range("D15").Style = $ws.Range("C15").Style
$ws.Range("E15").Value = '  +5.40%  '
$ws.Range("D16").Value = '1.812.64'
$ws.Range("E16").Value = '  +6.72%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001124'
$ws.Range("D17").Style = $ws.Range("C17").Style
$ws.Range("E17").Value = '  +4.87%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06791'
$ws.Range("D18").Style = $ws.Range("C18").Style
$ws.Range("E18").Value = '  +2.69%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '87.13'
$ws.Range("D19").Style = $ws.Range("C19").Style
$ws.Range("E19").Value = '  +7.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9991'
$ws.Range("D20").Style = $ws.Range("C20").Style
$ws.Range("E20").Value = '  +0.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.85'
$ws.Range("D21").Style = $ws.Range("C21").Style
$ws.Range("E21").Value = '  +10.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.547'
$ws.Range("D22").Style = $ws.Range("C22").Style
$ws.Range("E22").Value = '  +7.88%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.21'
$ws.Range("D23").Style = $ws.Range("C23").Style
$ws.Range("E23").Value = '  +1.73%  '
$ws.Range("D24").Value = '27.420.46'
$ws.Range("E24").Value = '  +6.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.474'
$ws.Range("D25").Style = $ws.Range("C25").Style
$ws.Range("E25").Value = '  +0.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.716'
$ws.Range("D26").Style = $ws.Range("C26").Style
$ws.Range("E26").Value = '  +9.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.49'
$ws.Range("D27").Style = $ws.Range("C27").Style
$ws.Range("E27").Value = '  +17.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.494'
$ws.Range("D28").Style = $ws.Range("C28").Style
$ws.Range("E28").Value = '  +14.27%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '154.58'
$ws.Range("D29").Style = $ws.Range("C29").Style
$ws.Range("E29").Value = '  +3.23%  '
$ws.Range("D30").Value = '2.016.23'
$ws.Range("E30").Value = '  +6.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '136.76'
$ws.Range("D31").Style = $ws.Range("C31").Style
$ws.Range("E31").Value = '  +6.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.387'
$ws.Range("D32").Style = $ws.Range("C32").Style
$ws.Range("E32").Value = '  +7.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.108'
$ws.Range("D33").Style = $ws.Range("C33").Style
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '13.83'
$ws.Range("D34").Style = $ws.Range("C34").Style
$ws.Range("E34").Value = '  +7.27%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08830'
$ws.Range("D35").Style = $ws.Range("C35").Style
$ws.Range("E35").Value = '  +4.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.724'
$ws.Range("D36").Style = $ws.Range("C36").Style
$ws.Range("E36").Value = '  +2.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.645'
$ws.Range("D37").Style = $ws.Range("C37").Style
$ws.Range("E37").Value = '  +6.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.7074'
$ws.Range("D38").Style = $ws.Range("C38").Style
$ws.Range("E38").Value = '  +15.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06565'
$ws.Range("D39").Style = $ws.Range("C39").Style
$ws.Range("E39").Value = '  +5.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2272'
$ws.Range("D40").Style = $ws.Range("C40").Style
$ws.Range("E40").Value = '  +7.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02419'
$ws.Range("D41").Style = $ws.Range("C41").Style
$ws.Range("E41").Value = '  +7.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.044'
$ws.Range("D42").Style = $ws.Range("C42").Style
$ws.Range("E42").Value = '  +5.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.259'
$ws.Range("D43").Style = $ws.Range("C43").Style
$ws.Range("E43").Value = '  -0.92%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.05'
$ws.Range("D44").Style = $ws.Range("C44").Style
$ws.Range("E44").Value = '  +3.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6609'
$ws.Range("D45").Style = $ws.Range("C45").Style
$ws.Range("E45").Value = '  +13.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9994'
$ws.Range("D46").Style = $ws.Range("C46").Style
$ws.Range("E46").Value = '  +0.05%  '
$ws.Range("E47").Value = '  +3.48%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.190'
$ws.Range("D48").Style = $ws.Range("C48").Style
$ws.Range("E48").Value = '  +9.23%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '132.93'
$ws.Range("D49").Style = $ws.Range("C49").Style
$ws.Range("E49").Value = '  +5.18%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07359'
$ws.Range("D50").Style = $ws.Range("C50").Style
$ws.Range("E50").Value = '  +2.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '80.73'
$ws.Range("D51").Style = $ws.Range("C51").Style
$ws.Range("E51").Value = '  +5.54%  '
